# Add model run v31 (2023_TM160_IPA_31) to the RTP2025 model run log.
# This inserts a new row right after the v30 row (row 42) and before the
# 2025 FBP Plus row (old row 43, which shifts down to row 44), shifting
# all subsequent rows down by one, then fills in the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 43, pushing existing rows 43+ down to 44+.
# Excel copies the format of the row above (row 42) into the new row,
# matching the "Past year" run-log entries' styling.
$ws.Rows("43:43").Insert()

$ws.Range("A43").Value = 2023
$ws.Range("B43").Value = "2023_TM160_IPA_31"
$ws.Range("C43").Value = "RTP2025_IP"
$ws.Range("D43").Value = "Past year"
$ws.Range("E43").Value = "Higher tolls, WFH remains at ~27%"
$ws.Range("F43").Value = "petrale"
$ws.Range("G43").Value = "n/a"
$ws.Range("H43").Value = "current"
$ws.Range("I43").Value = "BlueprintNetworks_v11\net_2023_Blueprint"
$ws.Range("J43").Value = "model2-b"
$ws.Range("K43").Value = "https://app.asana.com/0/1204085012544660/1205893933741809/f"
$ws.Range("L43").Value = 17.77
$ws.Range("M43").Value = "na"
$ws.Range("N43").Value = "na"
$ws.Range("O43").Value = 0.99
$ws.Range("P43").Value = 0.89
$ws.Range("Q43").Value = 120
$ws.Range("R43").Value = 0
$ws.Range("S43").Value = 45

# Match the author's final selection in the saved workbook.
$null = $ws.Range("E43").Select()
